$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A6").Value = "פרץ"
$ws.Range("B6").Value = "דביר"
$ws.Range("C6").Value = "0545885537"

$ws.Range("M7").Select() | Out-Null

